# Factsheet text-edit pass (COMM request): the "No. of 990 Filers w/ Gov
# Grants" counts (and the Overall sheet's single count) were stored as
# numbers; they need to become plain text cells -- using a thousands
# separator for the 2532 "Total" figures -- and the County sheet is
# missing its Total row entirely, so add it to match the other
# breakdown sheets (Congressional District / Size / Subsector).

$wb = $excel.ActiveWorkbook

# Helper: force $range to hold $value as TEXT, regardless of whether
# $value looks like a number/currency/percentage to Excel's normal
# type-sniffing on assignment. We momentarily mark the cell as "Text"
# (@) so the assignment isn't reinterpreted as a number, then restore
# the cell's style to Normal so no stray number-format style lingers
# on the cell afterwards.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---- Overall ----
$ws = $wb.Worksheets.Item("Overall")
Set-TextValue $ws.Range("A2") "2,532"

# ---- County ----
$ws = $wb.Worksheets.Item("County")
Set-TextValue $ws.Range("B2") "66"
Set-TextValue $ws.Range("B3") "231"
Set-TextValue $ws.Range("B4") "93"
Set-TextValue $ws.Range("B5") "142"
Set-TextValue $ws.Range("B6") "38"
Set-TextValue $ws.Range("B7") "44"
Set-TextValue $ws.Range("B8") "280"
Set-TextValue $ws.Range("B9") "37"
Set-TextValue $ws.Range("B10") "118"
Set-TextValue $ws.Range("B11") "63"
Set-TextValue $ws.Range("B12") "213"
Set-TextValue $ws.Range("B13") "176"
Set-TextValue $ws.Range("B14") "203"
Set-TextValue $ws.Range("B15") "168"
Set-TextValue $ws.Range("B16") "183"
Set-TextValue $ws.Range("B17") "124"
Set-TextValue $ws.Range("B18") "23"
Set-TextValue $ws.Range("B19") "120"
Set-TextValue $ws.Range("B20") "36"
Set-TextValue $ws.Range("B21") "134"
Set-TextValue $ws.Range("B22") "40"

# County sheet had no Total row (unlike the other breakdown sheets) --
# add row 23 to match.
Set-TextValue $ws.Range("A23") "Total"
Set-TextValue $ws.Range("B23") "2,532"
Set-TextValue $ws.Range("C23") '$5,748,555,291'
Set-TextValue $ws.Range("D23") "8.71%"
Set-TextValue $ws.Range("E23") "-15.65%"
Set-TextValue $ws.Range("F23") "69.47%"

# ---- Congressional District ----
$ws = $wb.Worksheets.Item("Congressional District")
Set-TextValue $ws.Range("B2") "172"
Set-TextValue $ws.Range("B3") "231"
Set-TextValue $ws.Range("B4") "251"
Set-TextValue $ws.Range("B5") "275"
Set-TextValue $ws.Range("B6") "201"
Set-TextValue $ws.Range("B7") "196"
Set-TextValue $ws.Range("B8") "232"
Set-TextValue $ws.Range("B9") "228"
Set-TextValue $ws.Range("B10") "210"
Set-TextValue $ws.Range("B11") "277"
Set-TextValue $ws.Range("B12") "124"
Set-TextValue $ws.Range("B13") "135"
Set-TextValue $ws.Range("B14") "2,532"

# ---- Size ----
$ws = $wb.Worksheets.Item("Size")
Set-TextValue $ws.Range("B2") "679"
Set-TextValue $ws.Range("B3") "713"
Set-TextValue $ws.Range("B4") "375"
Set-TextValue $ws.Range("B5") "245"
Set-TextValue $ws.Range("B6") "338"
Set-TextValue $ws.Range("B7") "182"
Set-TextValue $ws.Range("B8") "2,532"

# ---- Subsector ----
$ws = $wb.Worksheets.Item("Subsector")
Set-TextValue $ws.Range("B2") "197"
Set-TextValue $ws.Range("B3") "324"
Set-TextValue $ws.Range("B4") "73"
Set-TextValue $ws.Range("B5") "306"
Set-TextValue $ws.Range("B6") "19"
Set-TextValue $ws.Range("B7") "839"
Set-TextValue $ws.Range("B8") "17"
Set-TextValue $ws.Range("B9") "1"
Set-TextValue $ws.Range("B10") "167"
Set-TextValue $ws.Range("B11") "64"
Set-TextValue $ws.Range("B12") "502"
Set-TextValue $ws.Range("B13") "23"
Set-TextValue $ws.Range("B14") "2,532"
